$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D is treated as text so numeric-looking values
# (e.g. "1.001", "0.9999") are stored as text, not converted to numbers,
# matching the inlineStr cell type used throughout the sheet.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.860.10"
$ws.Range("E2").Value = "  -1.77%  "
$ws.Range("D3").Value = "1.832.78"
$ws.Range("E3").Value = "  -1.79%  "
$ws.Range("D4").Value = "0.9999"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "244.55"
$ws.Range("E5").Value = "  +0.51%  "
$ws.Range("D6").Value = "0.6898"
$ws.Range("E6").Value = "  -1.59%  "
$ws.Range("D7").Value = "1.0000"
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("D8").Value = "0.07699"
$ws.Range("E8").Value = "  -2.36%  "
$ws.Range("D9").Value = "0.3047"
$ws.Range("E9").Value = "  -2.52%  "
$ws.Range("D10").Value = "23.40"
$ws.Range("E10").Value = "  -3.93%  "
$ws.Range("D11").Value = "0.07806"
$ws.Range("E11").Value = "  +0.28%  "
$ws.Range("D12").Value = "1.839.77"
$ws.Range("E12").Value = "  -1.47%  "
$ws.Range("D13").Value = "5.096"
$ws.Range("E13").Value = "  -0.97%  "
$ws.Range("D14").Value = "92.03"
$ws.Range("E14").Value = "  -0.45%  "
$ws.Range("D15").Value = "0.6821"
$ws.Range("E15").Value = "  -2.45%  "
$ws.Range("D16").Value = "6.441"
$ws.Range("E16").Value = "  -1.24%  "
$ws.Range("D17").Value = "0.000008304"
$ws.Range("E17").Value = "  -3.51%  "
$ws.Range("D18").Value = "28.864.46"
$ws.Range("E18").Value = "  -1.82%  "
$ws.Range("D19").Value = "242.38"
$ws.Range("E19").Value = "  -2.48%  "
$ws.Range("D20").Value = "2.073.47"
$ws.Range("E20").Value = "  -2.33%  "
$ws.Range("E21").Value = "  -2.40%  "
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").Value = "7.445"
$ws.Range("E23").Value = "  -1.78%  "
$ws.Range("D24").Value = "1.001"
$ws.Range("E24").Value = "  -0.09%  "
$ws.Range("D25").Value = "0.1475"
$ws.Range("E25").Value = "  -3.99%  "
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "158.47"
$ws.Range("E26").Value = "  -1.32%  "
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").Value = "8.791"
$ws.Range("E27").Value = "  -2.05%  "
$ws.Range("D28").Value = "18.25"
$ws.Range("E28").Value = "  -2.38%  "
$ws.Range("D29").Value = "1.542"
$ws.Range("E29").Value = "  -2.61%  "
$ws.Range("D30").Value = "4.221"
$ws.Range("E30").Value = "  -1.75%  "
$ws.Range("D31").Value = "4.159"
$ws.Range("E31").Value = "  -2.04%  "
$ws.Range("D32").Value = "1.197"
$ws.Range("E32").Value = "  -0.75%  "
$ws.Range("D33").Value = "0.05084"
$ws.Range("E33").Value = "  -3.16%  "
$ws.Range("D34").Value = "0.7772"
$ws.Range("E34").Value = "  +2.51%  "
$ws.Range("D35").Value = "1.848"
$ws.Range("E35").Value = "  -2.11%  "
$ws.Range("D36").Value = "1.141"
$ws.Range("E36").Value = "  -3.54%  "
$ws.Range("D37").Value = "2.693"
$ws.Range("E37").Value = "  -0.50%  "
$ws.Range("D38").Value = "0.01852"
$ws.Range("E38").Value = "  -0.68%  "
$ws.Range("D39").Value = "1.240.15"
$ws.Range("E39").Value = "  -2.79%  "
$ws.Range("D40").Value = "2.696"
$ws.Range("E40").Value = "  -1.92%  "
$ws.Range("D41").Value = "0.9502"
$ws.Range("E41").Value = "  +5.97%  "
$ws.Range("D42").Value = "108.57"
$ws.Range("E42").Value = "  -1.07%  "
$ws.Range("D43").Value = "5.882"
$ws.Range("E43").Value = "  -1.02%  "
$ws.Range("D44").Value = "0.9997"
$ws.Range("E44").Value = "  -0.14%  "
$ws.Range("D45").Value = "9.601"
$ws.Range("E45").Value = "  +0.07%  "
$ws.Range("D46").Value = "0.00000000122"
$ws.Range("E46").Value = "  +1.08%  "
$ws.Range("D47").Value = "1.976.15"
$ws.Range("E47").Value = "  -2.23%  "
$ws.Range("D48").Value = "0.5159"
$ws.Range("E48").Value = "  -0.30%  "
$ws.Range("D49").Value = "63.92"
$ws.Range("E49").Value = "  -9.07%  "
$ws.Range("D50").Value = "1.745"
$ws.Range("E50").Value = "  -2.62%  "
$ws.Range("D51").Value = "6.919"
$ws.Range("E51").Value = "  -1.61%  "

# Reset the style index on column D back to the default/general style
# so no stray s="n" attribute is introduced (matches original formatting).
$ws.Range("D2:D51").Style = "Normal"
